$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 403 - this shifts the existing rows
# 403..426 down to 404..427, matching the diff.
$ws.Rows.Item(403).Insert()

# Populate the newly inserted row 403 with the new record.
$ws.Range("A403").Value = 5
$ws.Range("B403").Value = "Macroferia Regional de Talca"
$ws.Range("C403").Value = "Maule"
$ws.Range("D403").Value = 45013
$ws.Range("E403").Value = 7
$ws.Range("F403").Value = 100112006
$ws.Range("G403").Value = "Repollo"
$ws.Range("H403").Value = "Crespo record"
$ws.Range("I403").Value = "Primera"
$ws.Range("J403").Value = 3000
$ws.Range("K403").Value = 1000
$ws.Range("L403").Value = 1000
$ws.Range("M403").Value = 1000
$ws.Range("N403").Value = "$/unidad"
$ws.Range("O403").Value = "Región del Maule"
$ws.Range("P403").Value = 1000
$ws.Range("Q403").Value = 1
$ws.Range("R403").Value = "Hortaliza"
